$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Transformations sheet: replace rows 2-8 and append rows 9-16 with the new
# flo/od-split transformation names & equations.
# ---------------------------------------------------------------------------
$wsT = $wb.Worksheets.Item("Transformations")

$transformations = @(
    @("rem_none_samples_flo", "hcat(sg_1_none.flo,sg_2_none.flo,sg_3_none.flo,sg_4_none.flo) .- mean(control_M9_KC.flo)"),
    @("rem_atc_samples_flo", "hcat(sg_1_atc.flo,sg_2_atc.flo,sg_3_atc.flo,sg_4_atc.flo) .- mean(control_M9_KC_atc.flo)"),
    @("rem_IPTG_samples_flo", "hcat(sg_1_iptg.flo,sg_2_iptg.flo,sg_3_iptg.flo,sg_4_iptg.flo) .- mean(control_M9_KC_IPTG.flo)"),
    @("rem_atc_IPTG_samples_flo", "hcat(sg_1_atc_iptg.flo,sg_2_atc_iptg.flo,sg_3_atc_iptg.flo,sg_4_atc_iptg.flo) .- mean(control_M9_KC_atc_IPTG.flo)"),
    @("rem_single_flo", "hcat(sg1.flo,sg2.flo,sg3.flo,sg4.flo) .- mean(control_M9_K.flo)"),
    @("rem_wt_flo", "wt.flo .- mean(control_M9_NOAB.flo)"),
    @("rem_none_samples_od", "hcat(sg_1_none.OD,sg_2_none.OD,sg_3_none.OD,sg_4_none.OD) .- mean(control_M9_KC.OD)"),
    @("rem_atc_samples_od", "hcat(sg_1_atc.OD,sg_2_atc.OD,sg_3_atc.OD,sg_4_atc.OD) .- mean(control_M9_KC_atc.OD)"),
    @("rem_IPTG_samples_od", "hcat(sg_1_iptg.OD,sg_2_iptg.OD,sg_3_iptg.OD,sg_4_iptg.OD) .- mean(control_M9_KC_IPTG.OD)"),
    @("rem_atc_IPTG_samples_od", "hcat(sg_1_atc_iptg.OD,sg_2_atc_iptg.OD,sg_3_atc_iptg.OD,sg_4_atc_iptg.OD) .- mean(control_M9_KC_atc_IPTG.OD)"),
    @("rem_single_od", "hcat(sg1.OD,sg2.OD,sg3.OD,sg4.OD) .- mean(control_M9_K.OD)"),
    @("rem_wt_od", "wt.OD .- mean(control_M9_NOAB.OD)"),
    @("combine_samples_od", "hcat(rem_none_samples_od,rem_atc_samples_od,rem_IPTG_samples_od,rem_atc_IPTG_samples_od,rem_wt_od,rem_single_od)"),
    @("combine_samples_flo", "hcat(rem_none_samples_flo,rem_atc_samples_flo,rem_IPTG_samples_flo,rem_atc_IPTG_samples_flo,rem_wt_flo,rem_single_flo)"),
    @("flo_od", "(combine_samples_flo ./ combine_samples_od)")
)

$row = 2
foreach ($pair in $transformations) {
    $wsT.Cells.Item($row, 1).Value = $pair[0]
    $wsT.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

$wsT.Range("A16").Select() | Out-Null

# ---------------------------------------------------------------------------
# Views sheet: add a new "test" view row.
# ---------------------------------------------------------------------------
$wsV = $wb.Worksheets.Item("Views")
$wsV.Cells.Item(3, 1).Value = "test"
$wsV.Cells.Item(3, 2).Value = "plate_01_time.flo,flo_od"

$wsV.Range("B4").Select() | Out-Null
